# Se agrega funcionalidad para Menú
# Append new WhatsApp chat rows (Fecha / Mensaje) to the "Chats" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("03-05-2022 08:39", "hola"),
    @("03-05-2022 09:22", "Hola"),
    @("03-05-2022 09:22", "1"),
    @("03-05-2022 09:22", "Hola"),
    @("03-05-2022 09:22", "Hola"),
    @("03-05-2022 09:23", "Hola"),
    @("03-05-2022 09:24", "Adios")
)

$startRow = 27

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $fecha = $rows[$i][0]
    $mensaje = $rows[$i][1]

    $ws.Cells.Item($r, 1).Value = $fecha

    # "1" needs to be forced to text (the source data is a chat message,
    # not a number) so it round-trips as a shared string like the rest of
    # column B; everything else is already non-numeric-looking text.
    if ($mensaje -eq "1") {
        $cell = $ws.Cells.Item($r, 2)
        $cell.NumberFormat = "@"
        $cell.Value = $mensaje
        $cell.Style = "Normal"
    } else {
        $ws.Cells.Item($r, 2).Value = $mensaje
    }
}
